$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("B1").Value = "Email"

# Set column widths to match the new two-column layout as closely as the
# host's pixel-quantized ColumnWidth allows
$ws.Columns.Item(1).ColumnWidth = 14.3333333333333
$ws.Columns.Item(2).ColumnWidth = 28

$email = "jagjit.singh21802@gmail.com"

# Add a real mailto hyperlink on the first email cell (this also defines the
# built-in "Hyperlink" cell style used for the rest of the column)
$ws.Hyperlinks.Add($ws.Cells.Item(2, 2), "mailto:" + $email, "", "", $email)

# Repeat the same address/formatting down the remaining name rows (no extra
# functional hyperlinks were registered for these in the source workbook)
for ($r = 3; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $email
    $cell.Style = "Hyperlink"
}

# Trailing formatted-but-empty cell right below the table
$ws.Cells.Item(8, 2).Style = "Hyperlink"

# Move the active selection back to A7
$ws.Range("A7").Select() | Out-Null
